$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A8").Value = "Volume 31   Number  34"
$ws.Range("C9").Value = "Report Covering the Week  8/19/2024  Through  8/25/2024"

$ws.Range("D14").Value = "'0"
$ws.Range("E14").Value = "'***.*"
$ws.Range("C15").Value = 1
$ws.Range("D15").Value = "'0"
$ws.Range("E15").Value = "'***.*"
$ws.Range("F15").Value = 1
$ws.Range("H15").Value = 0
$ws.Range("I15").Value = 7
$ws.Range("K15").Value = 16.666666666666
$ws.Range("L15").Value = -22.222222222222
$ws.Range("M15").Value = 75
$ws.Range("N15").Value = -30
$ws.Range("D16").Value = 6
$ws.Range("E16").Value = -50
$ws.Range("F16").Value = 7
$ws.Range("G16").Value = 18
$ws.Range("H16").Value = -61.111111111111
$ws.Range("I16").Value = 50
$ws.Range("J16").Value = 77
$ws.Range("K16").Value = -35.064935064935
$ws.Range("L16").Value = -13.793103448275
$ws.Range("M16").Value = 11.111111111111
$ws.Range("N16").Value = -87.684729064039
$ws.Range("C17").Value = 3
$ws.Range("D17").Value = 3
$ws.Range("E17").Value = 0
$ws.Range("F17").Value = 8
$ws.Range("G17").Value = 20
$ws.Range("H17").Value = -60
$ws.Range("I17").Value = 75
$ws.Range("J17").Value = 78
$ws.Range("K17").Value = -3.846153846153
$ws.Range("L17").Value = 0
$ws.Range("M17").Value = 87.5
$ws.Range("N17").Value = -29.245283018867
$ws.Range("C18").Value = 3
$ws.Range("D18").Value = 2
$ws.Range("E18").Value = 50
$ws.Range("F18").Value = 11
$ws.Range("G18").Value = 13
$ws.Range("H18").Value = -15.384615384615
$ws.Range("I18").Value = 57
$ws.Range("J18").Value = 77
$ws.Range("K18").Value = -25.974025974026
$ws.Range("L18").Value = -55.46875
$ws.Range("M18").Value = -12.307692307692
$ws.Range("N18").Value = -92.578125
$ws.Range("C19").Value = 14
$ws.Range("D19").Value = 18
$ws.Range("E19").Value = -22.222222222222
$ws.Range("F19").Value = 53
$ws.Range("G19").Value = 55
$ws.Range("H19").Value = -3.636363636363
$ws.Range("I19").Value = 379
$ws.Range("J19").Value = 471
$ws.Range("K19").Value = -19.532908704883
$ws.Range("L19").Value = -10.823529411764
$ws.Range("M19").Value = -19.532908704883
$ws.Range("N19").Value = -74.112021857923
$ws.Range("C20").Value = 1
$ws.Range("D20").Value = 3
$ws.Range("E20").Value = -66.666666666666
$ws.Range("G20").Value = 5
$ws.Range("H20").Value = -20
$ws.Range("I20").Value = 18
$ws.Range("J20").Value = 38
$ws.Range("K20").Value = -52.631578947368
$ws.Range("L20").Value = -52.631578947368
$ws.Range("M20").Value = 5.882352941176
$ws.Range("N20").Value = -95.631067961165
$ws.Range("C21").Value = 25
$ws.Range("D21").Value = 32
$ws.Range("E21").Value = -21.875
$ws.Range("F21").Value = 84
$ws.Range("G21").Value = 113
$ws.Range("H21").Value = -25.663716814159
$ws.Range("I21").Value = 587
$ws.Range("J21").Value = 748
$ws.Range("K21").Value = -21.524064171123
$ws.Range("L21").Value = -19.918144611186
$ws.Range("M21").Value = -8.566978193146
$ws.Range("N21").Value = -81.459254579911
$ws.Range("D22").Value = "'0"
$ws.Range("E22").Value = "'***.*"
$ws.Range("G22").Value = 8
$ws.Range("H22").Value = -87.5
$ws.Range("L22").Value = 23.809523809523
$ws.Range("M22").Value = 30
$ws.Range("C24").Value = 27
$ws.Range("D24").Value = 18
$ws.Range("E24").Value = 50
$ws.Range("F24").Value = 86
$ws.Range("G24").Value = 85
$ws.Range("H24").Value = 1.176470588235
$ws.Range("I24").Value = 643
$ws.Range("J24").Value = 701
$ws.Range("K24").Value = -8.273894436519
$ws.Range("L24").Value = -27.180067950169
$ws.Range("M24").Value = 58.374384236453
$ws.Range("C25").Value = 22
$ws.Range("D25").Value = 15
$ws.Range("E25").Value = 46.666666666666
$ws.Range("F25").Value = 67
$ws.Range("G25").Value = 82
$ws.Range("H25").Value = -18.292682926829
$ws.Range("I25").Value = 518
$ws.Range("J25").Value = 587
$ws.Range("K25").Value = -11.75468483816
$ws.Range("L25").Value = -30.188679245283
$ws.Range("C26").Value = 6
$ws.Range("E26").Value = 0
$ws.Range("F26").Value = 13
$ws.Range("G26").Value = 22
$ws.Range("H26").Value = -40.90909090909
$ws.Range("I26").Value = 130
$ws.Range("J26").Value = 173
$ws.Range("K26").Value = -24.855491329479
$ws.Range("L26").Value = -15.584415584415
$ws.Range("M26").Value = -17.197452229299
$ws.Range("C27").Value = 1
$ws.Range("D27").Value = "'0"
$ws.Range("E27").Value = "'***.*"
$ws.Range("F27").Value = 1
$ws.Range("G27").Value = 1
$ws.Range("H27").Value = 0
$ws.Range("I27").Value = 8
$ws.Range("K27").Value = -11.111111111111
$ws.Range("L27").Value = -27.272727272727
$ws.Range("C28").Value = "'0"
$ws.Range("D28").Value = 6
$ws.Range("E28").Value = -100
$ws.Range("G28").Value = 10
$ws.Range("H28").Value = -80
$ws.Range("J28").Value = 48
$ws.Range("K28").Value = -29.166666666666
$ws.Range("D31").Value = 1
$ws.Range("E31").Value = -100
$ws.Range("G31").Value = 2
$ws.Range("J31").Value = 7
$ws.Range("K31").Value = 42.857142857142

$excel.CutCopyMode = 0

$ws.Range("A33").Copy()
$ws.Range("D14").PasteSpecial(-4122)
$ws.Range("A33").Copy()
$ws.Range("E14").PasteSpecial(-4122)
$ws.Range("J33").Copy()
$ws.Range("C15").PasteSpecial(-4122)
$ws.Range("A33").Copy()
$ws.Range("D15").PasteSpecial(-4122)
$ws.Range("A33").Copy()
$ws.Range("E15").PasteSpecial(-4122)
$ws.Range("J33").Copy()
$ws.Range("F15").PasteSpecial(-4122)
$ws.Range("J33").Copy()
$ws.Range("C17").PasteSpecial(-4122)
$ws.Range("J33").Copy()
$ws.Range("C20").PasteSpecial(-4122)
$ws.Range("J33").Copy()
$ws.Range("D20").PasteSpecial(-4122)
$ws.Range("K33").Copy()
$ws.Range("E20").PasteSpecial(-4122)
$ws.Range("A33").Copy()
$ws.Range("D22").PasteSpecial(-4122)
$ws.Range("A33").Copy()
$ws.Range("E22").PasteSpecial(-4122)
$ws.Range("J33").Copy()
$ws.Range("C27").PasteSpecial(-4122)
$ws.Range("A33").Copy()
$ws.Range("D27").PasteSpecial(-4122)
$ws.Range("A33").Copy()
$ws.Range("E27").PasteSpecial(-4122)
$ws.Range("J33").Copy()
$ws.Range("F27").PasteSpecial(-4122)
$ws.Range("A33").Copy()
$ws.Range("C28").PasteSpecial(-4122)
$ws.Range("J33").Copy()
$ws.Range("D31").PasteSpecial(-4122)
$ws.Range("K33").Copy()
$ws.Range("E31").PasteSpecial(-4122)
$excel.CutCopyMode = 0
